$wb = $excel.ActiveWorkbook

$rushing = $wb.Worksheets.Item("Rushing")
$receiving = $wb.Worksheets.Item("Receiving")

# --- Rushing sheet updates (Week 16 rushing stats) ---
$rushing.Range("C2").Value = 6
$rushing.Range("E2").Value = 5

$rushing.Range("C6").Value = 42
$rushing.Range("D6").Value = 29
$rushing.Range("E6").Value = 6
$rushing.Range("F6").Value = 13

$rushing.Range("C8").Value = 80
$rushing.Range("D8").Value = 55
$rushing.Range("E8").Value = 21
$rushing.Range("F8").Value = 29

$rushing.Range("C12").Value = 4

# --- Receiving sheet updates (Week 16 receiving stats) ---
$receiving.Range("C2").Value = 15
$receiving.Range("D2").Value = 11
$receiving.Range("G2").Value = 6
$receiving.Range("H2").Value = 5

$receiving.Range("C4").Value = 53
$receiving.Range("D4").Value = 47
$receiving.Range("G4").Value = 2
$receiving.Range("H4").Value = 2

$receiving.Range("C7").Value = 32
$receiving.Range("D7").Value = 18

$receiving.Range("C9").Value = 46
$receiving.Range("D9").Value = 33

$receiving.Range("C10").Value = 13
$receiving.Range("E10").Value = 8
$receiving.Range("G10").Value = 2

$receiving.Range("C11").Value = 7

$receiving.Range("C14").Value = 19
$receiving.Range("D14").Value = 12
$receiving.Range("E14").Value = 4
$receiving.Range("F14").Value = 3
$receiving.Range("G14").Value = 2
$receiving.Range("H14").Value = 2

# --- Active sheet moves back to Rushing ---
$rushing.Activate()
